$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.076.97"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.557.68"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9997"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3864"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3244"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.48"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -6.21%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07364"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.706"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.814"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "1.556.42"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001119"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06610"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.405"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9989"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").Value = "22.077.64"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.332"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.556"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.56"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.91"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.872"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "1.730.20"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.97"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.115"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.868"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.695"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -14.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08191"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.263"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06239"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02303"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.230"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2107"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.220"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.89"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5958"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.55"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.719"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5758"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.931"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06891"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.73%  "
